$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "67.619.49"
$ws.Range("E2").Value = "  -7.69%  "
$ws.Range("D3").Value = "3.674.39"
$ws.Range("E3").Value = "  -7.75%  "
$ws.Range("D4").Value = "'1.00"
$ws.Range("E4").Value = "  +0.09%  "
$ws.Range("D5").Value = "'567.05"
$ws.Range("E5").Value = "  -6.82%  "
$ws.Range("D6").Value = "'172.64"
$ws.Range("E6").Value = "  +0.47%  "
$ws.Range("D7").Value = "3.661.64"
$ws.Range("E7").Value = "  -7.93%  "
$ws.Range("D8").Value = "'0.620"
$ws.Range("E8").Value = "  -10.03%  "
$ws.Range("E9").Value = "  +0.10%  "
$ws.Range("D10").Value = "'0.699"
$ws.Range("E10").Value = "  -12.54%  "
$ws.Range("E11").Value = "  -13.52%  "
$ws.Range("D12").Value = "'51.11"
$ws.Range("E12").Value = "  -10.80%  "
$ws.Range("E13").Value = "  -13.98%  "
$ws.Range("D14").Value = "'10.42"
$ws.Range("E14").Value = "  -11.19%  "
$ws.Range("D15").Value = "4.265.47"
$ws.Range("E15").Value = "  -7.69%  "
$ws.Range("D16").Value = "3.670.81"
$ws.Range("E16").Value = "  -7.82%  "
$ws.Range("D17").Value = "'19.26"
$ws.Range("E17").Value = "  -8.38%  "
$ws.Range("E18").Value = "  -3.39%  "
$ws.Range("D19").Value = "'12.77"
$ws.Range("E19").Value = "  -10.96%  "
$ws.Range("E20").Value = "  -11.26%  "
$ws.Range("D21").Value = "67.345.14"
$ws.Range("E21").Value = "  -7.94%  "
$ws.Range("D22").Value = "'403.55"
$ws.Range("E22").Value = "  -13.14%  "
$ws.Range("D23").Value = "'4.40"
$ws.Range("E23").Value = "  -8.67%  "
$ws.Range("D24").Value = "'87.23"
$ws.Range("E24").Value = "  -9.87%  "
$ws.Range("E25").Value = "  -11.56%  "
$ws.Range("D26").Value = "'12.67"
$ws.Range("E26").Value = "  -11.50%  "
$ws.Range("D27").Value = "'10.62"
$ws.Range("E27").Value = "  -5.92%  "
$ws.Range("D28").Value = "'5.98"
$ws.Range("E28").Value = "  +1.48%  "
$ws.Range("E29").Value = "  -12.98%  "
$ws.Range("D30").Value = "'9.38"
$ws.Range("E30").Value = "  -12.11%  "
$ws.Range("D31").Value = "'32.39"
$ws.Range("E31").Value = "  -11.15%  "
$ws.Range("D32").Value = "'7.58"
$ws.Range("E32").Value = "  -6.46%  "
$ws.Range("D33").Value = "'12.41"
$ws.Range("E33").Value = "  -12.13%  "
$ws.Range("E34").Value = "  -11.65%  "
$ws.Range("D35").Value = "'64.57"
$ws.Range("E35").Value = "  -8.58%  "
$ws.Range("D36").Value = "'42.72"
$ws.Range("E36").Value = "  -14.68%  "
$ws.Range("D37").Value = "0.0₃0892"
$ws.Range("E37").Value = "  -13.38%  "
$ws.Range("D38").Value = "'579.41"
$ws.Range("E38").Value = "  -9.63%  "
$ws.Range("D39").Value = "'1.00"
$ws.Range("E39").Value = "  +0.00%  "
$ws.Range("D40").Value = "'0.393"
$ws.Range("E40").Value = "  -9.71%  "
$ws.Range("D41").Value = "'1.00"
$ws.Range("E41").Value = "  -0.10%  "
$ws.Range("D42").Value = "'0.133"
$ws.Range("E42").Value = "  -10.84%  "
$ws.Range("D43").Value = "'2.97"
$ws.Range("E43").Value = "  -8.97%  "
$ws.Range("E44").Value = "  -13.59%  "
$ws.Range("E45").Value = "  -11.25%  "
$ws.Range("D46").Value = "'2.55"
$ws.Range("E46").Value = "  -3.51%  "
$ws.Range("D47").Value = "'9.10"
$ws.Range("E47").Value = "  -14.14%  "
$ws.Range("E48").Value = "  -11.27%  "

$ws.Range("B49").Value = "WEMIXToken"
$ws.Range("C49").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D49").Value = "'2.68"
$ws.Range("E49").Value = "  -10.13%  "

$ws.Range("B50").Value = "ApeXProtocol"
$ws.Range("C50").Value = "https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex"
$ws.Range("D50").Value = "'3.15"
$ws.Range("E50").Value = "  -8.37%  "

$ws.Range("D51").Value = "2.699.65"
$ws.Range("E51").Value = "  -4.23%  "
